$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44830
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 12000
$ws.Range("P2").Value = 800

# Row 3
$ws.Range("D3").Value = 44771
$ws.Range("J3").Value = 40

# Row 4
$ws.Range("D4").Value = 44525
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 8000
$ws.Range("L4").Value = 8000
$ws.Range("M4").Value = 8000
$ws.Range("P4").Value = 533

# Row 5
$ws.Range("D5").Value = 44838
$ws.Range("J5").Value = 10

# Row 6
$ws.Range("D6").Value = 44749
$ws.Range("J6").Value = 50

# Row 7
$ws.Range("D7").Value = 44811
$ws.Range("J7").Value = 30
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 20000
$ws.Range("P7").Value = 1333

# Row 8
$ws.Range("D8").Value = 44757
$ws.Range("J8").Value = 30

# Row 9
$ws.Range("D9").Value = 44824
$ws.Range("J9").Value = 20

# Row 10
$ws.Range("D10").Value = 44819
$ws.Range("J10").Value = 100

# Row 11
$ws.Range("D11").Value = 44827
$ws.Range("J11").Value = 20

# Row 12
$ws.Range("D12").Value = 44826
$ws.Range("J12").Value = 50
$ws.Range("K12").Value = 20000
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = 20000
$ws.Range("P12").Value = 1333

# Row 13
$ws.Range("D13").Value = 44841
$ws.Range("J13").Value = 20
$ws.Range("K13").Value = 16000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 16000
$ws.Range("P13").Value = 1067

# Row 14
$ws.Range("D14").Value = 44508
$ws.Range("J14").Value = 40
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 10000
$ws.Range("P14").Value = 667

# Row 15
$ws.Range("D15").Value = 44839
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 16000
$ws.Range("L15").Value = 16000
$ws.Range("M15").Value = 16000
$ws.Range("P15").Value = 1067

# Row 16
$ws.Range("D16").Value = 44776
$ws.Range("J16").Value = 80

# Row 17
$ws.Range("D17").Value = 44845
$ws.Range("K17").Value = 16000
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 16000
$ws.Range("P17").Value = 1067

# Row 19
$ws.Range("D19").Value = 44812
$ws.Range("J19").Value = 80

# Row 20
$ws.Range("D20").Value = 44767
$ws.Range("J20").Value = 50

# Row 21
$ws.Range("D21").Value = 44825
$ws.Range("J21").Value = 30
$ws.Range("K21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = 20000
$ws.Range("P21").Value = 1333

# Row 22
$ws.Range("D22").Value = 44755
$ws.Range("J22").Value = 50

# Row 23
$ws.Range("D23").Value = 44813
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = 20000
$ws.Range("P23").Value = 1333

# Row 24
$ws.Range("D24").Value = 44756
$ws.Range("J24").Value = 80

# Row 25
$ws.Range("D25").Value = 44769
$ws.Range("J25").Value = 50
$ws.Range("K25").Value = 20000
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = 20000
$ws.Range("P25").Value = 1333

# Row 26
$ws.Range("D26").Value = 44837
$ws.Range("J26").Value = 80
$ws.Range("K26").Value = 16000
$ws.Range("L26").Value = 16000
$ws.Range("M26").Value = 16000
$ws.Range("P26").Value = 1067
